$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '26.553.18'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.77%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.669.15'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.24%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.007'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.37%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '219.63'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.59%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.5130'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.32%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.006'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.39%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.06439'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.28%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.2562'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.55%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '19.96'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.42%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07650'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.343'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.31%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.897.65'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.28%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.664.85'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.03%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.5572'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.87%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0₅8007'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.25%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '64.59'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.42%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '26.565.65'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.66%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.007'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.46%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '210.23'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.82%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.440'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.99%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.09'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.21%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.887'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.53%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.008'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.44%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '142.80'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.35%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.721'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.56%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.1168'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.92%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.989'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.66%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '15.67'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.92%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.21%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.20%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.348'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.63%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.192'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -6.72%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.576'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.75%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.13%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.376'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.85%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9231'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.57%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.5771'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.45%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.153.61'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +10.86%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.01584'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.007'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.41%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.8318'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.04%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.642'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.77%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '100.09'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.91%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.808.58'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.17%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0₈110'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.13%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '55.50'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.50%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.006'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.15%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.905'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.26%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.05135'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.00%  '
$cell.Style = 'Normal'

